$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowNum = 54

# Text columns (Date, Time, Weekday, Week) need to stay as literal text,
# not get auto-converted by Excel into a date serial / time serial / number.
$ws.Cells.Item($rowNum, 1).NumberFormat = "@"
$ws.Cells.Item($rowNum, 1).Value = "2025-02-06"
$ws.Cells.Item($rowNum, 1).Style = "Normal"

$ws.Cells.Item($rowNum, 2).NumberFormat = "@"
$ws.Cells.Item($rowNum, 2).Value = "09:05:01"
$ws.Cells.Item($rowNum, 2).Style = "Normal"

$ws.Cells.Item($rowNum, 3).NumberFormat = "@"
$ws.Cells.Item($rowNum, 3).Value = "Thursday"
$ws.Cells.Item($rowNum, 3).Style = "Normal"

$ws.Cells.Item($rowNum, 4).NumberFormat = "@"
$ws.Cells.Item($rowNum, 4).Value = "05"
$ws.Cells.Item($rowNum, 4).Style = "Normal"

# Numeric columns (Beijing ... Wuhan)
$numericValues = @(
    125807,
    141720,
    167089,
    157819,
    -1,
    142431,
    -1,
    -1,
    191122,
    115167,
    44693,
    28221,
    63069,
    -1,
    39331,
    -1
)

for ($i = 0; $i -lt $numericValues.Length; $i++) {
    $col = 5 + $i
    $ws.Cells.Item($rowNum, $col).Value = $numericValues[$i]
}
